# userCfg.xlsx maintenance edit
# 1) "emre" tab: a duplicate RAM/HyperX Fury row had snuck in - delete the dupe row.
# 2) "jan" tab: rebuilt with the current (non-working PC) component list and
#    moved to the end of the tab strip (it keeps getting in the way up front).

$wb = $excel.ActiveWorkbook

# --- 1) emre: remove the duplicate row (old row 4 == old row 3) ---
$emre = $wb.Worksheets.Item("emre")
$emre.Rows.Item(4).Delete() | Out-Null
$emre.Activate() | Out-Null
$emre.Range("A8").Select() | Out-Null

# --- 2) jan: drop the old sheet and re-add it at the end with fresh data ---
$jan = $wb.Worksheets.Item("jan")
$jan.Delete() | Out-Null

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newJan = $wb.Worksheets.Add($null, $lastSheet)
$newJan.Name = "jan"

$janData = @(
    @("groupComponent", "nameComponent"),
    @("Motherboard", "GA-Z87-HD3"),
    @("RAM", "CML8GX3M2A1600C9"),
    @("GPU", "GTX 980"),
    @("CPU", "i7 4770K"),
    @("PSU", "GS800"),
    @("Drive", "950 EVO")
)

for ($i = 0; $i -lt $janData.Count; $i++) {
    $row = $i + 1
    $newJan.Cells.Item($row, 1).Value = $janData[$i][0]
    $newJan.Cells.Item($row, 2).Value = $janData[$i][1]
}

# the "notworking" tab becomes the active tab after the reshuffle
$notworking = $wb.Worksheets.Item("notworking")
$notworking.Activate() | Out-Null
